$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.97547246983389, 50.05126917284872]"
$ws.Range("T2").Value = "[50.0140417861474, 50.06754030898669]"
$ws.Range("L3").Value = "[49.95707575026163, 50.04954033471086]"
$ws.Range("T3").Value = "[49.992841563223514, 50.05091364531852]"
